# Fruta / hortaliza, semanal
# Two new weekly observations are inserted at rows 16-17 (pushing the
# existing data for those dates down by two rows, down to what becomes
# row 131). Net effect: the sheet grows from 128 data rows (2-129) to
# 130 data rows (2-131).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 16, shifting all the
# existing data (old rows 16-129) down to rows 18-131.
$ws.Range("A16:A17").EntireRow.Insert()

# --- New row 16 ---
$ws.Cells.Item(16, 1).Value  = 3
$ws.Cells.Item(16, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(16, 3).Value  = "Coquimbo"
$ws.Cells.Item(16, 4).Value  = 44490
$ws.Cells.Item(16, 5).Value  = 5
$ws.Cells.Item(16, 6).Value  = "Fruta"
$ws.Cells.Item(16, 7).Value  = 100101
$ws.Cells.Item(16, 8).Value  = "Berries"
$ws.Cells.Item(16, 9).Value  = 100112025
$ws.Cells.Item(16, 10).Value = "Frutilla"
$ws.Cells.Item(16, 11).Value = "Sin especificar"
$ws.Cells.Item(16, 12).Value = "Especial"
$ws.Cells.Item(16, 13).Value = 98
$ws.Cells.Item(16, 14).Value = 7000
$ws.Cells.Item(16, 15).Value = 7000
$ws.Cells.Item(16, 16).Value = 7000
$ws.Cells.Item(16, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(16, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(16, 19).Value = 1000
$ws.Cells.Item(16, 20).Value = 7

# --- New row 17 ---
$ws.Cells.Item(17, 1).Value  = 3
$ws.Cells.Item(17, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(17, 3).Value  = "Coquimbo"
$ws.Cells.Item(17, 4).Value  = 44490
$ws.Cells.Item(17, 5).Value  = 5
$ws.Cells.Item(17, 6).Value  = "Fruta"
$ws.Cells.Item(17, 7).Value  = 100101
$ws.Cells.Item(17, 8).Value  = "Berries"
$ws.Cells.Item(17, 9).Value  = 100112025
$ws.Cells.Item(17, 10).Value = "Frutilla"
$ws.Cells.Item(17, 11).Value = "Sin especificar"
$ws.Cells.Item(17, 12).Value = "Segunda"
$ws.Cells.Item(17, 13).Value = 87
$ws.Cells.Item(17, 14).Value = 5000
$ws.Cells.Item(17, 15).Value = 5000
$ws.Cells.Item(17, 16).Value = 5000
$ws.Cells.Item(17, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(17, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(17, 19).Value = 714
$ws.Cells.Item(17, 20).Value = 7

# Make sure the D column (date) on the two new rows keeps the same
# date number format used by every other row in that column.
$ws.Range("D16:D17").NumberFormat = $ws.Range("D18").NumberFormat
